$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 2.12
$ws.Range("H2").Value = 2.67
$ws.Range("I2").Value = 4.15
$ws.Range("J2").Value = 2.85
$ws.Range("L2").Value = 4.85
$ws.Range("O2").Value = 1.6
$ws.Range("P2").Value = 2.05
$ws.Range("Q2").Value = 2.72
$ws.Range("R2").Value = 1.35
$ws.Range("U2").Value = 2.2
$ws.Range("V2").Value = 1.52
$ws.Range("W2").Value = 5
$ws.Range("X2").Value = 8.5
$ws.Range("Y2").Value = 9.5
$ws.Range("Z2").Value = 21
$ws.Range("AA2").Value = 23
$ws.Range("AB2").Value = 50
$ws.Range("AC2").Value = 4.45
$ws.Range("AD2").Value = 5.6
$ws.Range("AE2").Value = 21
$ws.Range("AF2").Value = 150
$ws.Range("AH2").Value = 7.8
$ws.Range("AI2").Value = 21
$ws.Range("AJ2").Value = 15.5
$ws.Range("AK2").Value = 80
$ws.Range("AN2").Value = 3.7
$ws.Range("AO2").Value = 11.75
$ws.Range("AP2").Value = 26
$ws.Range("AQ2").Value = 55
$ws.Range("AU2").Value = 8.25
$ws.Range("AW2").Value = 5.6
$ws.Range("AX2").Value = 27

# Row 5 updates
$ws.Range("P5").Value = 4.45
$ws.Range("U5").Value = 1.84
$ws.Range("V5").Value = 1.92
